# Insert a new weekly record for "Vega Monumental Concepción - Plátano" at
# row 540, pushing the previously existing rows 540:578 down to 541:579.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 540-578 down by one row.
$ws.Rows("540").Insert()

# Populate the newly inserted row 540 with the new record.
$ws.Range("A540").Value = 11
$ws.Range("B540").Value = "Vega Monumental Concepción"
$ws.Range("C540").Value = "Bíobío"
$ws.Range("D540").Value = 44826
$ws.Range("E540").Value = 8
$ws.Range("F540").Value = "Fruta"
$ws.Range("G540").Value = 100108
$ws.Range("H540").Value = "Tropicales y subtropicales"
$ws.Range("I540").Value = 100108006
$ws.Range("J540").Value = "Plátano"
$ws.Range("K540").Value = "Sin especificar"
$ws.Range("L540").Value = "Pintón"
$ws.Range("M540").Value = 1050
$ws.Range("N540").Value = 21000
$ws.Range("O540").Value = 22000
$ws.Range("P540").Value = 21476
$ws.Range("Q540").Value = "$/caja 20 kilos"
$ws.Range("R540").Value = "Ecuador"
$ws.Range("S540").Value = 1074
$ws.Range("T540").Value = 20
